$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# The "CARRO" (car) column is being renamed to "SAIDA" (departure), and its
# per-row values "CARRO 1".."CARRO 7" become "SAÍDA 1".."SAÍDA 7". Write the
# "SAÍDA n" values first (rows 2-20) so they claim the shared-string slots
# vacated by "CARRO 1".."CARRO 7", then set the new header text in B1 last -
# this reproduces the shared-string table order from the target file, where
# "STS" ends up right after "HUB" and the new "SAIDA" header string is
# appended at the very end.
$ws.Range("B2").Value = "SAÍDA 1"
$ws.Range("B3").Value = "SAÍDA 1"
$ws.Range("B4").Value = "SAÍDA 1"
$ws.Range("B5").Value = "SAÍDA 1"
$ws.Range("B6").Value = "SAÍDA 2"
$ws.Range("B7").Value = "SAÍDA 2"
$ws.Range("B8").Value = "SAÍDA 2"
$ws.Range("B9").Value = "SAÍDA 2"
$ws.Range("B10").Value = "SAÍDA 3"
$ws.Range("B11").Value = "SAÍDA 3"
$ws.Range("B12").Value = "SAÍDA 3"
$ws.Range("B13").Value = "SAÍDA 3"
$ws.Range("B14").Value = "SAÍDA 4"
$ws.Range("B15").Value = "SAÍDA 4"
$ws.Range("B16").Value = "SAÍDA 4"
$ws.Range("B17").Value = "SAÍDA 4"
$ws.Range("B18").Value = "SAÍDA 5"
$ws.Range("B19").Value = "SAÍDA 6"
$ws.Range("B20").Value = "SAÍDA 7"
$ws.Range("B1").Value = "SAIDA"

# Column B narrows slightly now that it holds "SAÍDA n" / "SAIDA" instead of
# "CARRO n" / "CARRO".
$ws.Columns("B").ColumnWidth = 7.67

# The active selection moved to H9.
$ws.Range("H9").Select()
